# Append the three newest admin-log transaction rows (7-9) that were
# recorded after the last export: sumit made three more deposits against
# account 12344.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2025-10-21 20:40:04", "sumit", "Deposit", "12344", 120),
    @("2025-10-21 20:40:12", "sumit", "Deposit", "12344", 77),
    @("2025-10-21 20:40:28", "sumit", "Deposit", "12344", 90)
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]                # Timestamp
    $ws.Cells.Item($r, 2).Value = $row[1]                # Admin
    $ws.Cells.Item($r, 3).Value = $row[2]                # Action
    # Leading apostrophe forces the account number to stay text (matches
    # the existing "12344" entries in column D, which are text too).
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]          # Account Affected
    $ws.Cells.Item($r, 5).Value = $row[4]                # Amount (numeric)
}
